$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L header (mirrors C1 "arena_height")
$ws.Range("L1").Value = "arena_height"

# Row 2 (XL_SPT)
$ws.Range("I2").Value = 3.68
$ws.Range("J2").Value = "13.077, 18.696"
$ws.Range("K2").Value = 0.378

# Row 3 (L_SPT)
$ws.Range("I3").Value = 1.865
$ws.Range("J3").Value = "6.455, 9.559"
$ws.Range("K3").Value = 0.22

# Row 4 (M_SPT)
$ws.Range("I4").Value = 0.87
$ws.Range("J4").Value = "3.32, 5.013"
$ws.Range("K4").Value = 0.121

# Row 5 (S_SPT)
$ws.Range("I5").Value = 0.532
$ws.Range("J5").Value = "2.65, 3.36"
$ws.Range("K5").Value = 0.074

# Update selection to match the recorded cursor position in the diff
$ws.Range("O12").Select()
